# Refresh the scraped coin quotes (Price / Volume(1h) columns) for the rows
# whose market data moved since the previous GitHub Actions run.
#
# The sheet stores these as literal text, not numbers/percentages. A plain
# Range.Value assignment of a numeric-looking string (e.g. "289.98" or
# "-3.64%") would make Excel silently reinterpret it as a Double (or a
# percentage) and rewrite the cell's number format. Prefixing the value with
# a leading apostrophe forces Excel to keep it as literal text (same as a
# user typing '289.98 into the Formula Bar); ClearFormats() afterwards
# strips the "quote prefix" marker Excel stamps on the cell's style for that
# apostrophe trick, so only the cell's text changes and its (default) style
# is left exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "289.98"
    "E2"  = "-3.64%"
    "D3"  = "30.48"
    "E3"  = "-5.61%"
    "D4"  = "4.945"
    "E4"  = "-0.30%"
    "D5"  = "0.07213"
    "E5"  = "-5.28%"
    "D6"  = "1.828"
    "E6"  = "-6.08%"
    "D7"  = "7.681"
    "E7"  = "-1.89%"
    "E8"  = "-0.55%"
    "D9"  = "0.8967"
    "E9"  = "-2.20%"
    "D10" = "0.1661"
    "E10" = "-5.90%"
    "D11" = "0.07737"
    "E11" = "-0.98%"
    "D12" = "0.08022"
    "E12" = "-5.92%"
    "D13" = "0.03034"
    "E13" = "-4.13%"
    "D14" = "0.1001"
    "E14" = "0.10%"
    "D15" = "0.001499"
    "E15" = "-1.04%"
    "D16" = "0.005776"
    "E16" = "-1.41%"
    "D18" = "3.470"
    "E18" = "0.24%"
    "D19" = "2.078"
    "E19" = "-3.50%"
    "D20" = "0.3319"
    "E20" = "-0.80%"
    "E21" = "-1.57%"
    "D22" = "4.042"
    "E22" = "-5.26%"
    "D23" = "0.2390"
    "E23" = "20.03%"
    "D24" = "0.04503"
    "E24" = "-0.21%"
    "D25" = "0.001217"
    "E25" = "-0.36%"
    "D26" = "0.004620"
    "E26" = "5.25%"
    "D27" = "0.0001302"
    "E27" = "4.10%"
    "D39" = "0.01567"
    "E39" = "-8.17%"
    "D40" = "0.04369"
    "E40" = "-6.61%"
    "D41" = "0.007330"
    "E41" = "-1.77%"
    "D42" = "0.009742"
    "D43" = "0.1301"
    "E43" = "-3.69%"
    "D44" = "0.002064"
    "E44" = "-11.51%"
    "D45" = "0.009523"
    "E45" = "-8.89%"
    "D46" = "0.00005969"
    "E46" = "-4.67%"
    "E47" = "0.05%"
    "D48" = "2.310"
    "E48" = "180.54%"
    "E50" = "0.05%"
    "E51" = "0.05%"
}

foreach ($address in $updates.Keys) {
    $rng = $ws.Range($address)
    $rng.Value = "'" + $updates[$address]
    $rng.ClearFormats()
}
